# Update cached market price / profit figures for the Ultima Profits sheets.
# (scheduled runner refresh - chore: update Sheets via scheduled runner)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1933.7778
$ws.Range("J40").Value = 1933.7778
$ws.Range("L40").Value = 1933.7778
$ws.Range("N40").Value = -2283.7778
# Row 64
$ws.Range("H64").Value = 2201218.8
$ws.Range("I64").Value = 4051797.8
$ws.Range("K64").Value = 4051797.8
$ws.Range("M64").Value = -4051549.8
# Row 67
$ws.Range("H67").Value = 2201218.8
$ws.Range("I67").Value = 4051797.8
$ws.Range("K67").Value = 4051797.8
$ws.Range("M67").Value = -4050939.8
# Row 74
$ws.Range("H74").Value = 4370.15
$ws.Range("I74").Value = 3674.875
$ws.Range("J74").Value = 4833.6665
$ws.Range("K74").Value = 3674.875
$ws.Range("L74").Value = 4833.6665
$ws.Range("M74").Value = -2738.875
$ws.Range("N74").Value = -6705.6665
# Row 76
$ws.Range("H76").Value = 3215.7673
$ws.Range("I76").Value = 3144.1177
$ws.Range("K76").Value = 3144.1177
$ws.Range("M76").Value = -2829.1177
# Row 77
$ws.Range("H77").Value = 4370.15
$ws.Range("I77").Value = 3674.875
$ws.Range("J77").Value = 4833.6665
$ws.Range("K77").Value = 18374.375
$ws.Range("L77").Value = 24168.3325
$ws.Range("M77").Value = -13694.375
$ws.Range("N77").Value = -33528.3325
# Row 79
$ws.Range("H79").Value = 3215.7673
$ws.Range("I79").Value = 3144.1177
$ws.Range("K79").Value = 3144.1177
$ws.Range("M79").Value = -2052.1177
# Row 138
$ws.Range("H138").Value = 1806.2703
$ws.Range("I138").Value = 915.2414
$ws.Range("J138").Value = 5036.25
$ws.Range("K138").Value = 2745.7242
$ws.Range("L138").Value = 15108.75
$ws.Range("M138").Value = 2394.2758
$ws.Range("N138").Value = -25388.75

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 14415.987
$ws.Range("I32").Value = 12387.687
$ws.Range("K32").Value = 12387.687
$ws.Range("M32").Value = -12100.687
# Row 117
$ws.Range("H117").Value = 30373.75
$ws.Range("J117").Value = 30373.75
$ws.Range("L117").Value = 30373.75
$ws.Range("N117").Value = -39551.75
# Row 122
$ws.Range("H122").Value = 8193
$ws.Range("I122").Value = 10638.286
$ws.Range("J122").Value = 2487.3333
$ws.Range("K122").Value = 31914.858
$ws.Range("L122").Value = 7461.999899999999
$ws.Range("M122").Value = -29464.858
$ws.Range("N122").Value = -12361.9999

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1939.2963
$ws.Range("I20").Value = 1977.5
$ws.Range("J20").Value = 1862.8889
$ws.Range("K20").Value = 1977.5
$ws.Range("L20").Value = 1862.8889
$ws.Range("M20").Value = -1730.5
$ws.Range("N20").Value = -2356.8889
# Row 80
$ws.Range("H80").Value = 25781.25
$ws.Range("J80").Value = 1466.6666
$ws.Range("L80").Value = 1466.6666
$ws.Range("N80").Value = -3462.6666
# Row 83
$ws.Range("H83").Value = 25781.25
$ws.Range("J83").Value = 1466.6666
$ws.Range("L83").Value = 7333.333000000001
$ws.Range("N83").Value = -17317.333
# Row 96
$ws.Range("H96").Value = 14712.333
$ws.Range("I96").Value = 5785.6
$ws.Range("J96").Value = 25870.75
$ws.Range("K96").Value = 5785.6
$ws.Range("L96").Value = 25870.75
$ws.Range("M96").Value = -3039.6
$ws.Range("N96").Value = -31362.75
# Row 105
$ws.Range("H105").Value = 4169.161
$ws.Range("I105").Value = 2620
$ws.Range("J105").Value = 4467.077
$ws.Range("K105").Value = 2620
$ws.Range("L105").Value = 4467.077
$ws.Range("M105").Value = -873
$ws.Range("N105").Value = -7961.077
# Row 118
$ws.Range("H118").Value = 7829.5454
$ws.Range("J118").Value = 7829.5454
$ws.Range("L118").Value = 7829.5454
$ws.Range("N118").Value = -11143.5454

$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 3682.7058
$ws.Range("I86").Value = 3422.889
$ws.Range("J86").Value = 3975
$ws.Range("K86").Value = 3422.889
$ws.Range("L86").Value = 3975
$ws.Range("M86").Value = -2299.889
$ws.Range("N86").Value = -6221
# Row 89
$ws.Range("H89").Value = 3682.7058
$ws.Range("I89").Value = 3422.889
$ws.Range("J89").Value = 3975
$ws.Range("K89").Value = 17114.445
$ws.Range("L89").Value = 19875
$ws.Range("M89").Value = -11498.445
$ws.Range("N89").Value = -31107

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 714.57574
$ws.Range("I5").Value = 285.27274
$ws.Range("J5").Value = 1573.1818
$ws.Range("K5").Value = 855.81822
$ws.Range("L5").Value = 4719.5454
$ws.Range("M5").Value = -743.81822
$ws.Range("N5").Value = -4943.5454
# Row 7
$ws.Range("H7").Value = 224.83333
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 224.83333
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 674.49999
$ws.Range("N7").Value = -898.49999
$ws.Range("M7").ClearContents()
# Row 11
$ws.Range("H11").Value = 118.35294
$ws.Range("I11").Value = 121
$ws.Range("J11").Value = 114.57143
$ws.Range("K11").Value = 363
$ws.Range("L11").Value = 343.71429
$ws.Range("M11").Value = -223
$ws.Range("N11").Value = -623.71429
# Row 44
$ws.Range("H44").Value = 1472
$ws.Range("J44").Value = 1628.2667
$ws.Range("L44").Value = 4884.800099999999
$ws.Range("N44").Value = -5680.800099999999
# Row 51
$ws.Range("H51").Value = 259.1111
$ws.Range("I51").Value = 259.1111
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 777.3333
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -317.3333
$ws.Range("N51").ClearContents()
# Row 69
$ws.Range("H69").Value = 2122.2222
$ws.Range("I69").Value = 2300
$ws.Range("J69").Value = 2100
$ws.Range("K69").Value = 6900
$ws.Range("L69").Value = 6300
$ws.Range("M69").Value = -6089
$ws.Range("N69").Value = -7922
# Row 72
$ws.Range("H72").Value = 2122.2222
$ws.Range("I72").Value = 2300
$ws.Range("J72").Value = 2100
$ws.Range("K72").Value = 20700
$ws.Range("L72").Value = 18900
$ws.Range("M72").Value = -16644
$ws.Range("N72").Value = -27012
# Row 118
$ws.Range("H118").Value = 2282.4167
$ws.Range("I118").Value = 2982.25
$ws.Range("J118").Value = 1932.5
$ws.Range("K118").Value = 8946.75
$ws.Range("L118").Value = 5797.5
$ws.Range("M118").Value = -7703.75
$ws.Range("N118").Value = -8283.5
# Row 131
$ws.Range("H131").Value = 1808.5555
$ws.Range("I131").Value = 2865.2942
$ws.Range("J131").Value = 1323.027
$ws.Range("K131").Value = 8595.882599999999
$ws.Range("L131").Value = 3969.081
$ws.Range("M131").Value = -3555.882599999999
$ws.Range("N131").Value = -14049.081
# Row 135
$ws.Range("H135").Value = 714.57574
$ws.Range("I135").Value = 285.27274
$ws.Range("J135").Value = 1573.1818
$ws.Range("K135").Value = 2567.45466
$ws.Range("L135").Value = 14158.6362
$ws.Range("M135").Value = -32.45465999999988
$ws.Range("N135").Value = -19228.6362

$ws = $wb.Worksheets.Item("GSM")
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
# Row 70
$ws.Range("H70").Value = 8700.584999999999
$ws.Range("J70").Value = 4005.5
$ws.Range("L70").Value = 4005.5
$ws.Range("N70").Value = -4545.5
# Row 73
$ws.Range("H73").Value = 8700.584999999999
$ws.Range("J73").Value = 4005.5
$ws.Range("L73").Value = 4005.5
$ws.Range("N73").Value = -5877.5
# Row 80
$ws.Range("H80").Value = 12304235
$ws.Range("I80").Value = 25643526
$ws.Range("J80").Value = 2103600.2
$ws.Range("K80").Value = 25643526
$ws.Range("L80").Value = 2103600.2
$ws.Range("M80").Value = -25642528
$ws.Range("N80").Value = -2105596.2
# Row 83
$ws.Range("H83").Value = 12304235
$ws.Range("I83").Value = 25643526
$ws.Range("J83").Value = 2103600.2
$ws.Range("K83").Value = 128217630
$ws.Range("L83").Value = 10518001
$ws.Range("M83").Value = -128212638
$ws.Range("N83").Value = -10527985
# Row 102
$ws.Range("H102").Value = 2915.7693
$ws.Range("I102").Value = 3647.158
$ws.Range("K102").Value = 3647.158
$ws.Range("M102").Value = -2025.158
# Row 113
$ws.Range("H113").Value = 46549.5
$ws.Range("I113").Value = 77718.92
$ws.Range("J113").Value = 1527
$ws.Range("K113").Value = 77718.92
$ws.Range("L113").Value = 1527
$ws.Range("M113").Value = -75548.92
$ws.Range("N113").Value = -5867
# Row 118
$ws.Range("H118").Value = 14805.556
$ws.Range("J118").Value = 14805.556
$ws.Range("L118").Value = 14805.556
$ws.Range("N118").Value = -18119.556

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 10005409
$ws.Range("I132").Value = 2522.606
$ws.Range("J132").Value = 29422776
$ws.Range("K132").Value = 7567.818000000001
$ws.Range("L132").Value = 88268328
$ws.Range("M132").Value = -5037.818000000001
$ws.Range("N132").Value = -88273388
